$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.022447764301858
$ws.Range("D2").Value = 1.043172239787303
$ws.Range("E2").Value = 1.02319583954603
$ws.Range("F2").Value = 1.047294593057895
$ws.Range("I2").Value = 1.03582502473026
$ws.Range("J2").Value = 1.027633547173591
$ws.Range("K2").Value = 1.045946855277735
$ws.Range("L2").Value = 1.026027990127646
$ws.Range("M2").Value = 1.050057631968986
$ws.Range("N2").Value = 1.013273405258344

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.023216197186943
$ws.Range("D3").Value = 1.043747512852678
$ws.Range("E3").Value = 1.023843306046228
$ws.Range("F3").Value = 1.048052146965428
$ws.Range("I3").Value = 1.035968399891887
$ws.Range("J3").Value = 1.028040973971687
$ws.Range("K3").Value = 1.046333726143478
$ws.Range("L3").Value = 1.026482718841767
$ws.Range("M3").Value = 1.05062714190919
$ws.Range("N3").Value = 1.01340944141458

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.023714194270559
$ws.Range("D4").Value = 1.044120306626345
$ws.Range("E4").Value = 1.02426331769611
$ws.Range("F4").Value = 1.048543407707096
$ws.Range("I4").Value = 1.036060312687169
$ws.Range("J4").Value = 1.028304703894755
$ws.Range("K4").Value = 1.046583907873848
$ws.Range("L4").Value = 1.026777330157236
$ws.Range("M4").Value = 1.050996059860621
$ws.Range("N4").Value = 1.013497464389435

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.023923734599559
$ws.Range("D5").Value = 1.044277159616812
$ws.Range("E5").Value = 1.024440141790286
$ws.Range("F5").Value = 1.048750188240265
$ws.Range("I5").Value = 1.036098746038854
$ws.Range("J5").Value = 1.028415597992068
$ws.Range("K5").Value = 1.046689046964206
$ws.Range("L5").Value = 1.026901272219351
$ws.Range("M5").Value = 1.051151248503688
$ws.Range("N5").Value = 1.013534468408046

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.023958927981172
$ws.Range("D6").Value = 1.04430350351621
$ws.Range("E6").Value = 1.024469846006471
$ws.Range("F6").Value = 1.048784922460025
$ws.Range("I6").Value = 1.036105187019037
$ws.Range("J6").Value = 1.028434218853052
$ws.Range("K6").Value = 1.046706698048727
$ws.Range("L6").Value = 1.026922087723577
$ws.Range("M6").Value = 1.051177310907157
$ws.Range("N6").Value = 1.013540681482447

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.023716993444532
$ws.Range("D7").Value = 1.04412240199354
$ws.Range("E7").Value = 1.024265679443771
$ws.Range("F7").Value = 1.048546169721504
$ws.Range("I7").Value = 1.036060827048999
$ws.Range("J7").Value = 1.028306185582167
$ws.Range("K7").Value = 1.046585312896001
$ws.Range("L7").Value = 1.02677898593529
$ws.Range("M7").Value = 1.050998133125844
$ws.Range("N7").Value = 1.013497958842844

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.022707299298007
$ws.Range("D8").Value = 1.043366540146654
$ws.Range("E8").Value = 1.023414433724283
$ws.Range("F8").Value = 1.047550388459885
$ws.Range("I8").Value = 1.035873656746664
$ws.Range("J8").Value = 1.027771217981333
$ws.Range("K8").Value = 1.0460776301718
$ws.Range("L8").Value = 1.026181590043386
$ws.Range("M8").Value = 1.050250015127326
$ws.Range("N8").Value = 1.013319379341018

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.020934067506165
$ws.Range("D9").Value = 1.042038945893415
$ws.Range("E9").Value = 1.021922616869433
$ws.Range("F9").Value = 1.045804011636977
$ws.Range("I9").Value = 1.035537282616241
$ws.Range("J9").Value = 1.026829341885899
$ws.Range("K9").Value = 1.045181948885795
$ws.Range("L9").Value = 1.025131811262482
$ws.Range("M9").Value = 1.048934935686079
$ws.Range("N9").Value = 1.013004708947402

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.019756045874095
$ws.Range("D10").Value = 1.041156917771171
$ws.Range("E10").Value = 1.020933695087288
$ws.Range("F10").Value = 1.044645489503367
$ws.Range("I10").Value = 1.035308671826877
$ws.Range("J10").Value = 1.026202044467368
$ws.Range("K10").Value = 1.044584191906493
$ws.Range("L10").Value = 1.024434001805472
$ws.Range("M10").Value = 1.048060480505593
$ws.Range("N10").Value = 1.012794963625111

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.019246953790837
$ws.Range("D11").Value = 1.040775735751112
$ws.Range("E11").Value = 1.020506839837869
$ws.Range("F11").Value = 1.044145224635236
$ws.Range("I11").Value = 1.03520865625389
$ws.Range("J11").Value = 1.02593058148056
$ws.Range("K11").Value = 1.044325224683881
$ws.Range("L11").Value = 1.024132345684769
$ws.Range("M11").Value = 1.047682392826893
$ws.Range("N11").Value = 1.012704155924385

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.01905800637016
$ws.Range("D12").Value = 1.040634261382385
$ws.Range("E12").Value = 1.020348492335587
$ws.Range("F12").Value = 1.043959614130527
$ws.Range("I12").Value = 1.035171352723506
$ws.Range("J12").Value = 1.025829773571727
$ws.Range("K12").Value = 1.044229014063792
$ws.Range("L12").Value = 1.024020373792512
$ws.Range("M12").Value = 1.047542039892516
$ws.Range("N12").Value = 1.012670428406514

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.019098529321557
$ws.Range("D13").Value = 1.040664602939647
$ws.Range("E13").Value = 1.020382449078081
$ws.Range("F13").Value = 1.043999418666287
$ws.Range("I13").Value = 1.035179361381552
$ws.Range("J13").Value = 1.025851396033316
$ws.Range("K13").Value = 1.044249652388155
$ws.Range("L13").Value = 1.024044388652024
$ws.Range("M13").Value = 1.047572142165576
$ws.Range("N13").Value = 1.012677662952338

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.019231332225122
$ws.Range("D14").Value = 1.040764039105828
$ws.Range("E14").Value = 1.020493746579027
$ws.Range("F14").Value = 1.044129877712425
$ws.Range("I14").Value = 1.035205575854996
$ws.Range("J14").Value = 1.025922248138717
$ws.Range("K14").Value = 1.044317272255806
$ws.Range("L14").Value = 1.024123088481454
$ws.Range("M14").Value = 1.047670789450753
$ws.Range("N14").Value = 1.012701367943291

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.019313176683016
$ws.Range("D15").Value = 1.040825320127154
$ws.Range("E15").Value = 1.020562347940598
$ws.Range("F15").Value = 1.044210285753449
$ws.Range("I15").Value = 1.035221707168068
$ws.Range("J15").Value = 1.025965905882867
$ws.Range("K15").Value = 1.044358932661136
$ws.Range("L15").Value = 1.024171588241986
$ws.Range("M15").Value = 1.047731580708589
$ws.Range("N15").Value = 1.012715973717243

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.019789853854342
$ws.Range("D16").Value = 1.04118223138115
$ws.Range("E16").Value = 1.02096205277601
$ws.Range("F16").Value = 1.044678719773801
$ws.Range("I16").Value = 1.035315287990441
$ws.Range("J16").Value = 1.026220064077961
$ws.Range("K16").Value = 1.044601376006199
$ws.Range("L16").Value = 1.024454032391304
$ws.Range("M16").Value = 1.048085584845642
$ws.Range("N16").Value = 1.012800990557557

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.020089129701709
$ws.Range("D17").Value = 1.041406312536199
$ws.Range("E17").Value = 1.021213141127133
$ws.Range("F17").Value = 1.04497292797276
$ws.Range("I17").Value = 1.035373714745229
$ws.Range("J17").Value = 1.026379534841005
$ws.Range("K17").Value = 1.044753419391405
$ws.Range("L17").Value = 1.024631337066425
$ws.Range("M17").Value = 1.048307792841853
$ws.Range("N17").Value = 1.012854323311877

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.020263788503833
$ws.Range("D18").Value = 1.041537086730119
$ws.Range("E18").Value = 1.021359727257165
$ws.Range("F18").Value = 1.045144667904006
$ws.Range("I18").Value = 1.03540769507492
$ws.Range("J18").Value = 1.026472566876276
$ws.Range("K18").Value = 1.044842090644241
$ws.Range("L18").Value = 1.024734804023499
$ws.Range("M18").Value = 1.048437456635988
$ws.Range("N18").Value = 1.012885432665934

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.020323358910487
$ws.Range("D19").Value = 1.041581689404517
$ws.Range("E19").Value = 1.021409731419038
$ws.Range("F19").Value = 1.045203249314174
$ws.Range("I19").Value = 1.035419264671734
$ws.Range("J19").Value = 1.026504290973395
$ws.Range("K19").Value = 1.044872322984131
$ws.Range("L19").Value = 1.024770091701541
$ws.Range("M19").Value = 1.048481677676003
$ws.Range("N19").Value = 1.012896040350632

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.02005701026455
$ws.Range("D20").Value = 1.04138226334344
$ws.Range("E20").Value = 1.021186188205829
$ws.Range("F20").Value = 1.044941348415827
$ws.Range("I20").Value = 1.035367456342542
$ws.Range("J20").Value = 1.026362423530851
$ws.Range("K20").Value = 1.044737107917029
$ws.Range("L20").Value = 1.024612308962594
$ws.Range("M20").Value = 1.048283946465594
$ws.Range("N20").Value = 1.012848601078127

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.019192220867159
$ws.Range("D21").Value = 1.0407347544801
$ws.Range("E21").Value = 1.020460966548508
$ws.Range("F21").Value = 1.044091454954889
$ws.Range("I21").Value = 1.035197860571241
$ws.Range("J21").Value = 1.025901383243969
$ws.Range("K21").Value = 1.044297360392959
$ws.Range("L21").Value = 1.024099911215995
$ws.Range("M21").Value = 1.047641737903572
$ws.Range("N21").Value = 1.012694387341186

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.018649374109348
$ws.Range("D22").Value = 1.040328298297231
$ws.Range("E22").Value = 1.020006181133578
$ws.Range("F22").Value = 1.04355831019537
$ws.Range("I22").Value = 1.035090342462325
$ws.Range("J22").Value = 1.02561165730246
$ws.Range("K22").Value = 1.044020766455426
$ws.Range("L22").Value = 1.023778190472196
$ws.Range("M22").Value = 1.047238452481093
$ws.Range("N22").Value = 1.012597441782005

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.018937063318877
$ws.Range("D23").Value = 1.040543705296614
$ws.Range("E23").Value = 1.020247157980928
$ws.Range("F23").Value = 1.043840824168381
$ws.Range("I23").Value = 1.035147423613488
$ws.Range("J23").Value = 1.025765231989167
$ws.Range("K23").Value = 1.044167403782491
$ws.Range("L23").Value = 1.02394869813186
$ws.Range("M23").Value = 1.047452193991448
$ws.Range("N23").Value = 1.012648832896437

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.020071523358458
$ws.Range("D24").Value = 1.041393129917995
$ws.Range("E24").Value = 1.02119836666772
$ws.Range("F24").Value = 1.044955617446835
$ws.Range("I24").Value = 1.035370284551351
$ws.Range("J24").Value = 1.026370155347846
$ws.Range("K24").Value = 1.044744478410809
$ws.Range("L24").Value = 1.024620906795121
$ws.Range("M24").Value = 1.04829472145113
$ws.Range("N24").Value = 1.012851186705419

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.021391770289711
$ws.Range("D25").Value = 1.042381634921402
$ws.Range("E25").Value = 1.022307304886865
$ws.Range("F25").Value = 1.046254492152966
$ws.Range("I25").Value = 1.035625015530307
$ws.Range("J25").Value = 1.027072735912365
$ws.Range("K25").Value = 1.045413622646297
$ws.Range("L25").Value = 1.025402850510764
$ws.Range("M25").Value = 1.049274524379563
$ws.Range("N25").Value = 1.013086054716482

